# Revert "Merging 0.1.8 w VitalSigns"
$wb = $excel.ActiveWorkbook

# --- Rename the two "Include ValueSet #N" sheets ---
$wb.Worksheets.Item("Include ValueSet #0").Name = "Include ValueSets"
$wb.Worksheets.Item("Include ValueSet #1").Name = "Include ValueSets 2"

# --- Update Metadata sheet values ---
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B3").Value = "0.1.6"
$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2023-05-05T10:50:04-05:00"
$ws.Range("B10").Value = "No display for ContactDetail"
$ws.Range("B11").Value = "No display for ContactDetail"

# Remove the old "Jurisdiction" row (row 12), shifting rows 13-16 up by one
$ws.Rows.Item(12).Delete()
